# Update column F (dSF) values on the active sheet to reflect the repulled
# data / recalculated means, as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    3  = 4
    4  = 9
    5  = -2
    6  = -1
    7  = 0
    10 = 4
    11 = 3
    12 = 1
    13 = -1
    14 = 2
    15 = 1
    16 = 2
    17 = -1
    19 = -2
    20 = -6
    21 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
